$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 6 data: SKU, Name, Landing Price, Selling Price, Units
$ws.Range("A6").Value = "0005"
$ws.Range("B6").Value = "Iron  Man Poster [Retro]"
$ws.Range("C6").Value = 140
$ws.Range("D6").Value = 249
$ws.Range("E6").Value = 1

# Match style of column A (text format) used by other SKU cells (row A2:A5)
$ws.Range("A6").NumberFormat = $ws.Range("A5").NumberFormat

# Update selection to mirror the diff (active cell moved to C6)
$ws.Range("C6").Select()
